$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying weekly data feed reshuffled which date/price row
# each record lands on. Columns A,B,C,E,F,G,O,R are identical for
# every data row, so the edit only needs to move D,H,I,J,K,L,M,N,P,Q
# values between rows 2-18 according to the new ordering.
$cols = @("D","H","I","J","K","L","M","N","P","Q")

# Snapshot current values before overwriting anything
$snapshot = @{}
for ($r = 2; $r -le 18; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping: target row -> row whose data should now appear there
$mapping = @{
    2 = 15
    3 = 6
    4 = 16
    5 = 13
    6 = 17
    7 = 18
    8 = 2
    9 = 8
    10 = 9
    11 = 7
    12 = 3
    13 = 10
    14 = 12
    15 = 5
    16 = 4
    17 = 14
    18 = 11
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $src = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $src[$c]
    }
}
